$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Que diseñará ... Microsoft Visual C# (con gestor de base de datos..."
#   -> split the run so " y Visual Basic" is inserted as its own run(s)
#      right after "Microsoft Visual C#".
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Microsoft Visual C#", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos1 = $rng1.End

$ins1 = $d.Range($insPos1, $insPos1)
$ins1.InsertAfter(" y Visual Basic")
$newRun1 = $d.Range($insPos1, $insPos1 + 15)
# Temporarily change the formatting so the engine keeps this as a distinct
# run instead of merging it back into its neighbours; restored at the end.
$newRun1.Font.Size = 10

# ---------------------------------------------------------------------------
# Edit 2: "... sistema, la cantidad de ( Q 1000" + "0.00 ) quetzales..."
#   -> "Q 10000.00" becomes "Q 15,000.00", written as 4 runs: "...Q 1", "5",
#      ",", "00", followed by the untouched "0.00 ) quetzales..." run.
# ---------------------------------------------------------------------------
$rngA = $d.Content
$rngA.Find.Execute("el cliente se obliga a pagar al proveedor por el desarrollo del", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# Also protect the preceding run from being re-absorbed by the edit below.
$rngA.Font.Size = 10

$rng2 = $d.Content
$rng2.Find.Execute(" sistema, la cantidad de ( Q 1000", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$delStart = $rng2.End - 3
$delEnd = $rng2.End
$delRange = $d.Range($delStart, $delEnd)
$delRange.Delete()

$insPosA = $delStart
$insA = $d.Range($insPosA, $insPosA)
$insA.InsertAfter("5")
$newRunA = $d.Range($insPosA, $insPosA + 1)
$newRunA.Font.Size = 9

$insPosB = $insPosA + 1
$insB = $d.Range($insPosB, $insPosB)
$insB.InsertAfter(",")
$newRunB = $d.Range($insPosB, $insPosB + 1)
$newRunB.Font.Size = 8

$insPosC = $insPosB + 1
$insC = $d.Range($insPosC, $insPosC)
$insC.InsertAfter("00")
$newRunC = $d.Range($insPosC, $insPosC + 2)
$newRunC.Font.Size = 7

# ---------------------------------------------------------------------------
# Restore the original font size (11 pt == sz 22) on every run we touched,
# now that all insertions/deletions are complete. Doing this last prevents
# the intermediate edits from re-merging the runs we just split apart.
# ---------------------------------------------------------------------------
$newRun1.Font.Size = 11
$rngA.Font.Size = 11
$newRunA.Font.Size = 11
$newRunB.Font.Size = 11
$newRunC.Font.Size = 11
